$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '306.83'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.52%'
$ws.Range("E2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '38.89'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '6.91%'
$ws.Range("E3").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.104'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.83%'
$ws.Range("E4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08082'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.64%'
$ws.Range("E5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.940'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-5.17%'
$ws.Range("E6").ClearFormats()
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.963'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.14%'
$ws.Range("E7").ClearFormats()
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9315'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.32%'
$ws.Range("E8").ClearFormats()
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1476'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '3.17%'
$ws.Range("E9").ClearFormats()
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1930'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.15%'
$ws.Range("E10").ClearFormats()
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09146'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.58%'
$ws.Range("E11").ClearFormats()
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03501'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.36%'
$ws.Range("E12").ClearFormats()
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09792'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.30%'
$ws.Range("E13").ClearFormats()
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001402'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.15%'
$ws.Range("E14").ClearFormats()
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005840'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-5.96%'
$ws.Range("E15").ClearFormats()
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.787'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.33%'
$ws.Range("E16").ClearFormats()
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.185'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.70%'
$ws.Range("E17").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.00%'
$ws.Range("E19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1302'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.80%'
$ws.Range("E20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.563'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-4.92%'
$ws.Range("E21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2413'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '3.21%'
$ws.Range("E22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04369'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.02%'
$ws.Range("E23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001236'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.45%'
$ws.Range("E24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004284'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-12.89%'
$ws.Range("E25").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.18%'
$ws.Range("E26").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '0.17%'
$ws.Range("E39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05071'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.42%'
$ws.Range("E40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007522'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.74%'
$ws.Range("E41").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '2.02%'
$ws.Range("E42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1351'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.89%'
$ws.Range("E43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002122'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.29%'
$ws.Range("E44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009093'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-6.51%'
$ws.Range("E45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006206'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-1.29%'
$ws.Range("E46").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.18%'
$ws.Range("E47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003105'
$ws.Range("D48").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.18%'
$ws.Range("E50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.18%'
$ws.Range("E51").ClearFormats()
